$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in previously empty Subject/Summary cells for existing rows 3-5
$ws.Range("B3").Value = "tesing 2 subject"
$ws.Range("C3").Value = "summary test 2"

$ws.Range("B4").Value = "subject test 3"
$ws.Range("C4").Value = "summary test 3"

$ws.Range("B5").Value = "4th subject test"
$ws.Range("C5").Value = "5th summary test"

# Add new row 6
$ws.Range("A6").Value = 25
$ws.Range("B6").Value = "test 1"
$ws.Range("C6").Value = "syummary test"
$ws.Range("D6").Value = "17/06/2025"

# Add new row 7
$ws.Range("A7").Value = 5
$ws.Range("B7").Value = "dhskfh"
$ws.Range("C7").Value = "fdov8ysdf"
$ws.Range("D7").Value = "17/06/2025"
